# Apply "assign trip excel data changes" to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Vehicle Registration Number values in column A
$ws.Range("A5").Value = "MP04CT1237"
$ws.Range("A6").Value = "MP04CT1237"
$ws.Range("A7").Value = "MP04CT1237"
$ws.Range("A14").Value = "MP04CT1239"
$ws.Range("A15").Value = "MP04CT1239"

# Update the active selection to match the saved view state (B16)
$ws.Range("B16").Select()
